$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheets
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
# (shared by Overview!E2, Overview!F2, zh-cn!C2, de-de!C2)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# ---------------------------------------------------------------------------
# Column widths
#   Overview E,F  : 17.2159881591797 -> 29.9777047293527
#   zh-cn/de-de C : 17.2159881591797 -> 29.9777047293527
#   zh-cn/de-de I : 18.6506053379604 -> 40
#   zh-cn/de-de J : 21.7054770333426 -> 40
# ---------------------------------------------------------------------------
$wideNarrow = 29.166666666666668   # quantizes to stored width 30 (closest achievable to 29.9777047293527)
$wideFull   = 39.166666666666664   # quantizes to stored width 40 exactly

$overview.Columns.Item(5).ColumnWidth = $wideNarrow   # E
$overview.Columns.Item(6).ColumnWidth = $wideNarrow   # F

foreach ($ws in @($zhcn, $dede)) {
    $ws.Columns.Item(3).ColumnWidth  = $wideNarrow  # C
    $ws.Columns.Item(9).ColumnWidth  = $wideFull     # I
    $ws.Columns.Item(10).ColumnWidth = $wideFull     # J
}

# ---------------------------------------------------------------------------
# zh-cn / de-de row 2: Latest Target File (I2) -> add hyperlink + text first,
# for both sheets, THEN fix up the visual style on both - this keeps the
# number of transient/unused style records the Hyperlinks.Add call leaves
# behind to a minimum.
# ---------------------------------------------------------------------------
$zhI2 = $zhcn.Range("I2")
$zhI2.Value = "855bebd0-e14f-4b76-afd8-bfd13c3e8764.md"
$zhcn.Hyperlinks.Add($zhI2, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c8a8df3dba688edf6ac190d4667facb86d47ce21/e2e/855bebd0-e14f-4b76-afd8-bfd13c3e8764.md", "", "", "855bebd0-e14f-4b76-afd8-bfd13c3e8764.md") | Out-Null

$deI2 = $dede.Range("I2")
$deI2.Value = "855bebd0-e14f-4b76-afd8-bfd13c3e8764.md"
$dede.Hyperlinks.Add($deI2, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c8a8df3dba688edf6ac190d4667facb86d47ce21/e2e/855bebd0-e14f-4b76-afd8-bfd13c3e8764.md", "", "", "855bebd0-e14f-4b76-afd8-bfd13c3e8764.md") | Out-Null

$zhI2.Style = "HyperLink"
$zhI2.Font.Underline = $true
$zhI2.Font.Color = 15570276   # FF6495ED (BGR-encoded) - matches existing hyperlink font

$deI2.Style = "HyperLink"
$deI2.Font.Underline = $true
$deI2.Font.Color = 15570276

# ---------------------------------------------------------------------------
# zh-cn row 2: Latest Handback File (J2) + Latest Handback DateTime (K2)
# ---------------------------------------------------------------------------
$zhcn.Range("J2").Value = "855bebd0-e14f-4b76-afd8-bfd13c3e8764.f96f0d11bed2aa4c5b25d9175ddaedaab9b4cc71.zh-cn.xlf"

# zh-cn's "Latest Handback DateTime" becomes 2016-08-26 04:56:33
$zhcn.Range("K2").Value = "2016-08-26 04:56:33"

# ---------------------------------------------------------------------------
# de-de row 2: Latest Handback File (J2) + Latest Handback DateTime (K2)
# ---------------------------------------------------------------------------
$dede.Range("J2").Value = "855bebd0-e14f-4b76-afd8-bfd13c3e8764.f96f0d11bed2aa4c5b25d9175ddaedaab9b4cc71.de-de.xlf"

$dede.Range("K2").Value = "2016-08-26 04:56:40"

Write-Host "Handback report generated."
